# The deck's single Design ("Integral") is swapped back to the stock
# "Office Theme" colour scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink),
# i.e. ppt/theme/theme2.xml's <a:clrScheme> values change from the
# Integral palette to the default Office palette. (ppt/theme/theme1.xml,
# the Notes Master's theme, keeps the Integral values it already holds.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index order (PpThemeColorSchemeIndex): 1 Dark1, 2 Light1, 3 Dark2,
# 4 Light2, 5 Accent1 .. 10 Accent6, 11 Hyperlink, 12 FollowedHyperlink.
# Values are the standard VBA RGB() encoding (R + G*256 + B*65536) for
# the stock "Office Theme" palette.
$officeColors = @(
    0,          # dk1     000000
    16777215,   # lt1     FFFFFF
    6968388,    # dk2     44546A
    15132391,   # lt2     E7E6E6
    13998939,   # accent1 5B9BD5
    3243501,    # accent2 ED7D31
    10855845,   # accent3 A5A5A5
    49407,      # accent4 FFC000
    12874308,   # accent5 4472C4
    4697456,    # accent6 70AD47
    12673797,   # hlink   0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
